# Add two new WI ("2025-26" site) hikes worth of readings: rows 986-997
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$site = "2025-26"

# Data for the new rows: row -> (B date-serial, C depth, D SWE-or-rain, E note)
# Column D / E are omitted (left blank) where the source row has no value.
$rows = @(
    @{ R = 986; B = 46029.5;              C = 1.83; D = $null;  E = "Extrapolated" },
    @{ R = 987; B = 46030.384722222225;   C = 1.82; D = $null;  E = $null },
    @{ R = 988; B = 46031.466666666667;   C = 1.83; D = $null;  E = $null },
    @{ R = 989; B = 46032.5;              C = 1.83; D = 0.03;   E = "SWE; Extrapolated" },
    @{ R = 990; B = 46033.520833333336;   C = 1.82; D = 0.04;   E = "SWE" },
    @{ R = 991; B = 46034.534722222219;   C = 1.83; D = $null;  E = $null },
    @{ R = 992; B = 46035.498611111114;   C = 1.82; D = 0.14000000000000001; E = "Rain (not SWE)" },
    @{ R = 993; B = 46036.520833333336;   C = 1.84; D = 0.01;   E = "SWE" },
    @{ R = 994; B = 46037.5;              C = 1.83; D = $null;  E = "Extrapolated" },
    @{ R = 995; B = 46038.583333333336;   C = 1.86; D = 0.4;    E = "SWE" },
    @{ R = 996; B = 46039.594444444447;   C = 1.84; D = 0.04;   E = "SWE" },
    @{ R = 997; B = 46040.602777777778;   C = 1.86; D = 0.01;   E = "SWE" }
)

foreach ($row in $rows) {
    $r = $row.R

    $ws.Cells.Item($r, 1).Value = $site

    $ws.Cells.Item($r, 2).Value = $row.B
    # Match the existing date/time number format used by column B (copy from
    # the row directly above, which already carries the right style index).
    $ws.Range("B" + ($r - 1)).Copy()
    $ws.Cells.Item($r, 2).PasteSpecial(-4122)

    $ws.Cells.Item($r, 3).Value = $row.C

    if ($null -ne $row.D) {
        $ws.Cells.Item($r, 4).Value = $row.D
    }

    if ($null -ne $row.E) {
        $ws.Cells.Item($r, 5).Value = $row.E
    }
}

$excel.CutCopyMode = $false

# Move the frozen-pane selection down to the new last row, matching how
# Excel leaves the active cell after appending rows at the bottom.
[void]$ws.Range("D998").Select()

Write-Host "Added rows 986-997"
